$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 9, shifting existing rows 9:105 down to 10:106
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data
$ws.Cells.Item(9, 1).Value = 11
$ws.Cells.Item(9, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(9, 3).Value = "Bíobío"
$ws.Cells.Item(9, 4).Value = 44630
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 100112032
$ws.Cells.Item(9, 7).Value = "Zapallo italiano"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 10000
$ws.Cells.Item(9, 12).Value = 11000
$ws.Cells.Item(9, 13).Value = 10500
$ws.Cells.Item(9, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 175
$ws.Cells.Item(9, 17).Value = 60
$ws.Cells.Item(9, 18).Value = "Hortaliza"
